$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14-41 down to 15-42.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with its data.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44536
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103001
$ws.Range("J14").Value = "Cereza"
$ws.Range("K14").Value = "Santina"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 400
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 11000
$ws.Range("P14").Value = 10500
$ws.Range("Q14").Value = "`$/bandeja 8 kilos"
$ws.Range("R14").Value = "Provincia de Curicó"
$ws.Range("S14").Value = 1312
$ws.Range("T14").Value = 8
